$d = $word.ActiveDocument

# The first paragraph is currently empty; the trailing paragraph carries
# the "_GoBack" bookmark. Move that bookmark onto the first paragraph
# (Bookmarks.Add re-targets an existing bookmark name, removing it from
# wherever it was), then add the text "n" to the first paragraph.
$r = $d.Paragraphs(1).Range
$d.Bookmarks.Add("_GoBack", $r)

$p1 = $d.Paragraphs(1).Range
$p1.InsertBefore("n")
